$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", reusing the same formatting as the other
# header cells (B1:G1), e.g. G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save column values (0/1) for rows 2-15 (H2..H15)
$saveValues = @(1, 0, 0, 1, 0, 1, 0, 0, 0, 0, 1, 1, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
